# "both folders are updated on 19 Nov"
# The current (until-now) Sheet1 becomes the dated snapshot "Sheet1_2(19Nov)"
# (kept in the first tab position, with a freshly added K-column of numbers),
# while an untouched copy of the original Sheet1 data keeps living on under
# the plain name "Sheet1" right after it. Sheet2 / Sheet3 just shift along.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")

# Duplicate Sheet1 and place the duplicate immediately before itself - the
# original keeps its name/position+1, the new copy lands in Sheet1's old slot.
$sheet1.Copy($sheet1)

# The copy is now the first sheet in the workbook; rename it to the dated tab.
$dated = $wb.Worksheets.Item(1)
$dated.Name = "Sheet1_2(19Nov)"

# The original sheet (still literally named "Sheet1") is now one slot later.
$original = $wb.Worksheets.Item("Sheet1")

# --- Fill in the 19-Nov K column ("units/day, 10 days") on the dated sheet ---
$dated.Activate()
$dated.Range("K8").Value = 1555
$dated.Range("K9").Value = 188
$dated.Range("K10").Value = 1868
$dated.Range("K11").Value = 7750
$dated.Range("K12").Value = 86
$dated.Range("K13").Value = 693
$dated.Range("K14").Value = 12
$dated.Range("K15").Value = 5880

# Reflect where the sheet was left scrolled/selected after the data entry.
$dated.Range("K15").Select()

# The plain "Sheet1" copy is left exactly as it was before this edit, just
# scrolled/selected the way the duplicate naturally starts out.
$original.Activate()
$original.Range("J16").Select()

$dated.Activate()
